$d = $word.ActiveDocument

function Set-ParagraphText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# Functional requirements bullet list starts at paragraph 4 (1-indexed):
#   4  Website must allow users to search for the poll they want
#   5  There must be a page where you can view the latest content
#   6  There must be an option to view how many times an item has been voted
#   7  Website must have a place to store the users data such as login details, and a range of poles.
#   8  You must full in all your details when creating your account, so full name, username, password, and email address.
#   9  To create a post you must fill in all relevant details
#   10 You must be logged in to create a post
#   11 You must choose a survey type ... in order to post the survey.
#   12 Password must contain letters and numbers
#   13 There should be a navigation bar where you can easily navigate the pages

Set-ParagraphText 4  "Website must have a place to store the users data such as login details, and a range of poles."
Set-ParagraphText 5  "You must full in all your details when creating your account, so full name, username, password, and email address."
Set-ParagraphText 6  "To create a post you must fill in all relevant details "
Set-ParagraphText 7  "You must choose a survey type such as text based, number based, colour based in order to post the survey."
Set-ParagraphText 8  "There should be a navigation bar where you can easily navigate the pages"
Set-ParagraphText 9  "You must fill in username to add a user "
Set-ParagraphText 10 "You must put in your details correctly to log in"
Set-ParagraphText 11 "Passwords must be not visible when inputted"
Set-ParagraphText 12 "Passwords must be hashed on the database so they are not visible to anyone making it encrypted"

# The old paragraph 13 (navigation bar bullet) is now redundant - its content
# moved up to paragraph 8 above - so the trailing paragraph is removed entirely.
$d.Paragraphs.Item(13).Range.Delete()
